$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, matching the formatting of H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new data columns I and J for rows 2-12
$data = @{
    2  = @(9, 9)
    3  = @(10, 10)
    4  = @(9, 9)
    5  = @(8, 8)
    6  = @(4, 4)
    7  = @(8, 9)
    8  = @(5, 6)
    9  = @(6, 6)
    10 = @(1, 1)
    11 = @(5, 5)
    12 = @(3, 3)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
